# Scheduled runner refresh of market-board derived leve-profit figures
# (currentAveragePrice* / LevePrice* / LeveProfit* columns) across the
# ALC / ARM / BSM / CRP / CUL / GSM / LTW / WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1860.6666
$ws.Range("I129").Value = 410
$ws.Range("J129").Value = 2116.6667
$ws.Range("K129").Value = 1230
$ws.Range("L129").Value = 6350.000100000001
$ws.Range("M129").Value = 3770
$ws.Range("N129").Value = -16350.0001

$ws.Range("H131").Value = 5867.6523
$ws.Range("I131").Value = 993.7143
$ws.Range("J131").Value = 8000
$ws.Range("K131").Value = 2981.1429
$ws.Range("L131").Value = 24000
$ws.Range("M131").Value = 2058.8571
$ws.Range("N131").Value = -34080

$ws.Range("H137").Value = 9091816
$ws.Range("I137").Value = 572.7
$ws.Range("J137").Value = 16667851
$ws.Range("K137").Value = 1718.1
$ws.Range("L137").Value = 50003553
$ws.Range("M137").Value = 831.8999999999999
$ws.Range("N137").Value = -50008653

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 50353.332
$ws.Range("J24").Value = 50353.332
$ws.Range("L24").Value = 50353.332
$ws.Range("N24").Value = -51101.332

$ws.Range("H32").Value = 10943.708
$ws.Range("I32").Value = 11589.741
$ws.Range("J32").Value = 9005.611000000001
$ws.Range("K32").Value = 11589.741
$ws.Range("L32").Value = 9005.611000000001
$ws.Range("M32").Value = -11302.741
$ws.Range("N32").Value = -9579.611000000001

$ws.Range("H74").Value = 9436030
$ws.Range("I74").Value = 10870969
$ws.Range("J74").Value = 6430.5713
$ws.Range("K74").Value = 10870969
$ws.Range("L74").Value = 6430.5713
$ws.Range("M74").Value = -10870095
$ws.Range("N74").Value = -8178.5713

$ws.Range("H77").Value = 9436030
$ws.Range("I77").Value = 10870969
$ws.Range("J77").Value = 6430.5713
$ws.Range("K77").Value = 54354845
$ws.Range("L77").Value = 32152.8565
$ws.Range("M77").Value = -54350477
$ws.Range("N77").Value = -40888.85649999999

$ws.Range("H100").Value = 50353.332
$ws.Range("J100").Value = 50353.332
$ws.Range("L100").Value = 50353.332
$ws.Range("N100").Value = -52517.332

$ws.Range("H101").Value = 50602
$ws.Range("J101").Value = 50602
$ws.Range("L101").Value = 50602
$ws.Range("N101").Value = -57092

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2090.4285
$ws.Range("I20").Value = 1938.0834
$ws.Range("J20").Value = 3004.5
$ws.Range("K20").Value = 1938.0834
$ws.Range("L20").Value = 3004.5
$ws.Range("M20").Value = -1691.0834
$ws.Range("N20").Value = -3498.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8500936
$ws.Range("I31").Value = 2443324.2
$ws.Range("J31").Value = 33337144
$ws.Range("K31").Value = 2443324.2
$ws.Range("L31").Value = 33337144
$ws.Range("M31").Value = -2443029.2
$ws.Range("N31").Value = -33337734

$ws.Range("H34").Value = 8500936
$ws.Range("I34").Value = 2443324.2
$ws.Range("J34").Value = 33337144
$ws.Range("K34").Value = 2443324.2
$ws.Range("L34").Value = 33337144
$ws.Range("M34").Value = -2443122.2
$ws.Range("N34").Value = -33337548

$ws.Range("H132").Value = 9435706
$ws.Range("I132").Value = 11365088
$ws.Range("J132").Value = 3171.111
$ws.Range("K132").Value = 34095264
$ws.Range("L132").Value = 9513.332999999999
$ws.Range("M132").Value = -34092734
$ws.Range("N132").Value = -14573.333

$ws.Range("H134").Value = 322819.8
$ws.Range("I134").Value = 1125
$ws.Range("J134").Value = 1134716.2
$ws.Range("K134").Value = 3375
$ws.Range("L134").Value = 3404148.6
$ws.Range("M134").Value = -840
$ws.Range("N134").Value = -3409218.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 221.8
$ws.Range("I11").Value = 186.66667
$ws.Range("J11").Value = 274.5
$ws.Range("K11").Value = 560.00001
$ws.Range("L11").Value = 823.5
$ws.Range("M11").Value = -420.00001
$ws.Range("N11").Value = -1103.5

$ws.Range("H26").Value = 254.54546
$ws.Range("I26").Value = 144.44444
$ws.Range("J26").Value = 750
$ws.Range("K26").Value = 433.33332
$ws.Range("L26").Value = 2250
$ws.Range("M26").Value = -145.33332
$ws.Range("N26").Value = -2826

$ws.Range("H68").Value = 712.6667
$ws.Range("I68").Value = 334
$ws.Range("J68").Value = 807.3333
$ws.Range("K68").Value = 1002
$ws.Range("L68").Value = 2421.9999
$ws.Range("M68").Value = -191
$ws.Range("N68").Value = -4043.9999

$ws.Range("H71").Value = 712.6667
$ws.Range("I71").Value = 334
$ws.Range("J71").Value = 807.3333
$ws.Range("K71").Value = 3006
$ws.Range("L71").Value = 7265.9997
$ws.Range("M71").Value = 1050
$ws.Range("N71").Value = -15377.9997

$ws.Range("H93").Value = 7280
$ws.Range("I93").Value = 1400
$ws.Range("K93").Value = 4200
$ws.Range("M93").Value = -2328

$ws.Range("H109").Value = 3896
$ws.Range("I109").Value = 1000
$ws.Range("J109").Value = 4620
$ws.Range("K109").Value = 3000
$ws.Range("L109").Value = 13860
$ws.Range("M109").Value = -1960
$ws.Range("N109").Value = -15940

$ws.Range("H115").Value = 5777.778
$ws.Range("I115").Value = 0
$ws.Range("J115").Value = 5777.778
$ws.Range("K115").Value = 0
$ws.Range("L115").Value = 17333.334
$ws.Range("M115").Value = ""
$ws.Range("N115").Value = -19683.334

$ws.Range("H118").Value = 1188.08
$ws.Range("J118").Value = 1247.1305
$ws.Range("L118").Value = 3741.3915
$ws.Range("N118").Value = -6227.3915

$ws.Range("H131").Value = 1653.6316
$ws.Range("I131").Value = 2578.3333
$ws.Range("J131").Value = 1226.8462
$ws.Range("K131").Value = 7734.999899999999
$ws.Range("L131").Value = 3680.5386
$ws.Range("M131").Value = -2694.999899999999
$ws.Range("N131").Value = -13760.5386

$ws.Range("H140").Value = 3117
$ws.Range("I140").Value = 1645
$ws.Range("J140").Value = 5570.3335
$ws.Range("K140").Value = 4935
$ws.Range("L140").Value = 16711.0005
$ws.Range("M140").Value = 245
$ws.Range("N140").Value = -27071.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2566031
$ws.Range("I122").Value = 4168299.5
$ws.Range("J122").Value = 2401.6
$ws.Range("K122").Value = 12504898.5
$ws.Range("L122").Value = 7204.799999999999
$ws.Range("M122").Value = -12502448.5
$ws.Range("N122").Value = -12104.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H31").Value = 1689
$ws.Range("I31").Value = 1323.3334
$ws.Range("J31").Value = 1871.8334
$ws.Range("K31").Value = 1323.3334
$ws.Range("L31").Value = 1871.8334
$ws.Range("M31").Value = -1075.3334
$ws.Range("N31").Value = -2367.8334

$ws.Range("H60").Value = 20061
$ws.Range("J60").Value = 20061
$ws.Range("L60").Value = 20061
$ws.Range("N60").Value = -21079

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7602.857
$ws.Range("I62").Value = 4215.4
$ws.Range("J62").Value = 10682.363
$ws.Range("K62").Value = 4215.4
$ws.Range("L62").Value = 10682.363
$ws.Range("M62").Value = -3591.4
$ws.Range("N62").Value = -11930.363

$ws.Range("H65").Value = 7602.857
$ws.Range("I65").Value = 4215.4
$ws.Range("J65").Value = 10682.363
$ws.Range("K65").Value = 21077
$ws.Range("L65").Value = 53411.815
$ws.Range("M65").Value = -17957
$ws.Range("N65").Value = -59651.815
